$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.907957333333333
$ws.Range("H2").Value = 5.723871999999999
$ws.Range("I2").Value = 0.02117870709996734
$ws.Range("J2").Value = 0.02117870709996734
$ws.Range("M2").Value = 0.5273163333333333
$ws.Range("N2").Value = 1.581949
$ws.Range("Q2").Value = 1.006097065169778
$ws.Range("R2").Value = 9.054873586527998
$ws.Range("S2").Value = 0.02117870709996734
$ws.Range("T2").Value = 0.02117870709996734

# Row 3
$ws.Range("I3").Value = 0.8627330973450646
$ws.Range("J3").Value = 0.8627330973450646
$ws.Range("M3").Value = 0.5273163333333333
$ws.Range("N3").Value = 1.581949
$ws.Range("Q3").Value = 40.98424106658711
$ws.Range("R3").Value = 368.858169599284
$ws.Range("S3").Value = 0.8627330973450646
$ws.Range("T3").Value = 0.8627330973450646

# Row 4
$ws.Range("G4").Value = 10.458208
$ws.Range("H4").Value = 31.374624
$ws.Range("I4").Value = 0.116088195554968
$ws.Range("J4").Value = 0.116088195554968
$ws.Range("M4").Value = 0.5273163333333333
$ws.Range("N4").Value = 1.581949
$ws.Range("Q4").Value = 5.514783895797333
$ws.Range("R4").Value = 49.633055062176
$ws.Range("S4").Value = 0.116088195554968
$ws.Range("T4").Value = 0.116088195554968
